# Updated symbol list on Fri Feb  3 23:51:20 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# rows whose underlying market data changed. Values are kept as text
# (matching the workbook's existing inline-string cell format) rather
# than letting Excel auto-convert them to numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $cell = $ws.Range($address)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$updates = [ordered]@{
    "D2"  = "332.02"; "E2"  = "2.08%"
    "D3"  = "41.17";  "E3"  = "3.01%"
    "D4"  = "5.733";  "E4"  = "-2.32%"
    "D5"  = "0.08196"; "E5" = "2.66%"
    "D6"  = "2.042";  "E6"  = "6.70%"
    "D7"  = "8.749"
    "E8"  = "-1.45%"
    "D9"  = "2.993";  "E9"  = "1.78%"
    "D10" = "0.9209"; "E10" = "-1.94%"
    "D11" = "0.1246"; "E11" = "-0.94%"
    "D12" = "0.1948"; "E12" = "-0.66%"
    "D13" = "8.308";  "E13" = "-5.82%"
    "D14" = "0.09428"; "E14" = "3.05%"
    "D15" = "0.03622"; "E15" = "1.44%"
    "E16" = "9.73%"
    "D17" = "0.001298"; "E17" = "-0.31%"
    "D18" = "0.006224"; "E18" = "0.88%"
    "D19" = "3.384";  "E19" = "1.14%"
    "E20" = "-1.17%"
    "E21" = "-1.24%"
    "D22" = "0.2649"; "E22" = "9.56%"
    "D23" = "0.04428"; "E23" = "-0.59%"
    "E24" = "-0.13%"
    "D25" = "0.004317"; "E25" = "-0.70%"
    "E26" = "8.51%"
    "D39" = "0.02767"; "E39" = "14.32%"
    "D40" = "0.05510"; "E40" = "4.87%"
    "D41" = "0.007621"; "E41" = "2.33%"
    "D42" = "0.009947"
    "D43" = "0.1422"; "E43" = "0.78%"
    "D44" = "0.002130"; "E44" = "1.26%"
    "D45" = "0.01180"; "E45" = "11.97%"
    "D46" = "0.00006739"; "E46" = "-1.93%"
    "E48" = "59.77%"
    "D49" = "0.003007"; "E49" = "4.35%"
}

foreach ($addr in $updates.Keys) {
    Set-TextValue $addr $updates[$addr]
}
